$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new requirement rows (req4 and req5), following the same pattern
# as the existing rows 3-5 (id, title, description), copying their style.

# Row 6: req4
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "req4"
$ws.Range("D6").Value = "Software requirement description four"

# Row 7: req5
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "req5"
$ws.Range("D7").Value = "Software requirement description five"

# Copy formatting from row 5 (the last populated data row) onto rows 6 and 7
# so they pick up the same borders/fills/number formats as the other rows.
$ws.Range("B5:D5").Copy() | Out-Null
$ws.Range("B6:D7").PasteSpecial(-4122) | Out-Null

# Re-apply the values since PasteSpecial(xlPasteFormats) should not touch
# values, but ensure they are still set correctly.
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "req4"
$ws.Range("D6").Value = "Software requirement description four"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "req5"
$ws.Range("D7").Value = "Software requirement description five"

$excel.CutCopyMode = 0
